# support header & footer
#
# The header/footer paragraphs each contain a run whose text ("header" /
# "footer") is wrapped in a <w:commentRangeStart>/<w:commentRangeEnd> pair
# that has no backing comment part. We need to:
#   1. drop the now-stray commentRangeStart/commentRangeEnd markers, and
#   2. rename the placeholder text to MDWORD-HEADER / MDWORD-FOOTER.
#
# Word's Range.Find can change the run text but cannot reach the
# (non-textual) comment-range markers, so we rebuild each paragraph via
# Range.InsertXML with the exact original run/paragraph formatting minus
# the comment markers.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

# ---- Header -----------------------------------------------------------
$hdr = $sec.Headers(1)
$hdrRange = $hdr.Range
$hdrXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/dummy.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="163AB4E7" w14:textId="1481A90A" w:rsidR="00870D41" w:rsidRPr="0002329F" w:rsidRDefault="00870D41"><w:pPr><w:pStyle w:val="ae"/><w:rPr><w:rFonts w:ascii="隶书" w:eastAsia="隶书"/></w:rPr></w:pPr><w:r w:rsidRPr="0002329F"><w:rPr><w:rFonts w:ascii="隶书" w:eastAsia="隶书" w:hint="eastAsia"/></w:rPr><w:t/></w:r><w:r w:rsidRPr="0002329F"><w:rPr><w:rFonts w:ascii="隶书" w:eastAsia="隶书" w:hint="eastAsia"/></w:rPr><w:t>MDWORD-HEADER</w:t></w:r><w:r w:rsidRPr="0002329F"><w:rPr><w:rFonts w:ascii="隶书" w:eastAsia="隶书" w:hint="eastAsia"/></w:rPr><w:t/></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>
'@
$hdrRange.InsertXML($hdrXml)

# ---- Footer -------------------------------------------------------------
$ftr = $sec.Footers(1)
$ftrRange = $ftr.Range
$ftrXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/dummy.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="33847476" w14:textId="42701653" w:rsidR="000E54C2" w:rsidRPr="00690A4D" w:rsidRDefault="00FA3820" w:rsidP="00FA3820"><w:pPr><w:pStyle w:val="af0"/><w:rPr><w:rFonts w:ascii="隶书" w:eastAsia="隶书" w:hAnsi="宋体"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="隶书" w:eastAsia="隶书" w:hAnsi="宋体"/></w:rPr><w:t xml:space="preserve">123 </w:t></w:r><w:r w:rsidR="000E54C2" w:rsidRPr="00690A4D"><w:rPr><w:rFonts w:ascii="隶书" w:eastAsia="隶书" w:hAnsi="宋体" w:hint="eastAsia"/></w:rPr><w:t/></w:r><w:r w:rsidR="000E54C2" w:rsidRPr="00690A4D"><w:rPr><w:rFonts w:ascii="隶书" w:eastAsia="隶书" w:hAnsi="宋体" w:hint="eastAsia"/></w:rPr><w:t>MDWORD-FOOTER</w:t></w:r><w:r w:rsidR="000E54C2" w:rsidRPr="00690A4D"><w:rPr><w:rFonts w:ascii="隶书" w:eastAsia="隶书" w:hAnsi="宋体" w:hint="eastAsia"/></w:rPr><w:t/></w:r><w:r><w:rPr><w:rFonts w:ascii="隶书" w:eastAsia="隶书" w:hAnsi="宋体"/></w:rPr><w:t xml:space="preserve"> 789</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>
'@
$ftrRange.InsertXML($ftrXml)

Write-Output "header/footer comment markers stripped; placeholders renamed"
